$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# ---------------------------------------------------------------------------
# 1) Extend Table1 by one row at the bottom (ref A8:K137 -> A8:K138).
# ---------------------------------------------------------------------------
$tbl.ListRows.Add()

# ---------------------------------------------------------------------------
# 2) Fill in the new leave entries on rows 88-91 (row numbers unaffected by
#    the row-92 insert below, since they sit above it). Done before the
#    insert so new shared strings are created in the same order Excel would
#    (row 89's strings before row 92's "2024" label).
# ---------------------------------------------------------------------------
$ws.Range("C88").Value = 1.25

$ws.Range("B89").Value = "VL(62-0-0)"
$ws.Range("C89").Value = 1.25
$ws.Range("D89").Value = 62
$ws.Range("K89").Value = "10/2 - 12/29/2023"

$ws.Range("C90").Value = 1.25

$ws.Range("C91").Value = 1.25

# ---------------------------------------------------------------------------
# 3) Insert a new blank row at row 92 (a "2024" year-separator row), which
#    shifts every row from 92 down through the newly-added 138 down by one.
# ---------------------------------------------------------------------------
$ws.Rows("92:92").Insert()

# Restore the normal body-row formatting on the freshly inserted row 92
# (Excel's Insert leaves it with "default" formatting otherwise).
$ws.Range("A91:K91").Copy()
$ws.Range("A92").PasteSpecial(-4122)
$ws.Range("G92").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# A92 becomes a "2024" text label like the other year-header rows
# (A10="2017", A14="2018", ... A79="2023"), styled the same way.
$ws.Range("A92").NumberFormat = "@"
$ws.Range("A92").Value = "2024"
$ws.Range("A79").Copy()
$ws.Range("A92").PasteSpecial(-4122)

# Fix up the calculated column formula on the new last row (138), which the
# table-extend step above leaves without a formula.
$ws.Range("G138").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# ---------------------------------------------------------------------------
# 4) Refresh the view state to match (pane/selection).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("F5").Select()
$excel.ActiveWindow.ScrollRow = 82
$ws.Range("F98").Select()
